$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 575
$ws.Range("I31").Value = 575
$ws.Range("K31").Value = 1725
$ws.Range("M31").Value = -1495
$ws.Range("H32").Value = 4317
$ws.Range("I32").Value = 3999.5
$ws.Range("J32").Value = 4475.75
$ws.Range("K32").Value = 3999.5
$ws.Range("L32").Value = 4475.75
$ws.Range("M32").Value = -3673.5
$ws.Range("N32").Value = -5127.75
$ws.Range("H106").Value = 4639.091
$ws.Range("I106").Value = 2859.2856
$ws.Range("J106").Value = 7753.75
$ws.Range("K106").Value = 2859.2856
$ws.Range("L106").Value = 7753.75
$ws.Range("M106").Value = -2228.2856
$ws.Range("N106").Value = -9015.75
$ws.Range("H113").Value = 1870.3914
$ws.Range("I113").Value = 1599.5
$ws.Range("J113").Value = 1896.1904
$ws.Range("K113").Value = 1599.5
$ws.Range("L113").Value = 1896.1904
$ws.Range("M113").Value = 1654.5
$ws.Range("N113").Value = -8404.190399999999
$ws.Range("H116").Value = 3080583
$ws.Range("J116").Value = 3928.5715
$ws.Range("L116").Value = 3928.5715
$ws.Range("N116").Value = -10812.5715
$ws.Range("H132").Value = 4559.7812
$ws.Range("I132").Value = 3833.8147
$ws.Range("J132").Value = 8480
$ws.Range("K132").Value = 11501.4441
$ws.Range("L132").Value = 25440
$ws.Range("M132").Value = -8971.444100000001
$ws.Range("N132").Value = -30500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1228.2727
$ws.Range("I45").Value = 930.5
$ws.Range("J45").Value = 1398.4286
$ws.Range("K45").Value = 930.5
$ws.Range("L45").Value = 1398.4286
$ws.Range("M45").Value = -553.5
$ws.Range("N45").Value = -2152.4286
$ws.Range("H61").Value = 1886.875
$ws.Range("I61").Value = 1277.85
$ws.Range("K61").Value = 1277.85
$ws.Range("M61").Value = -1065.85
$ws.Range("H132").Value = 2020.9474
$ws.Range("I132").Value = 1607.5834
$ws.Range("J132").Value = 2729.5715
$ws.Range("K132").Value = 4822.7502
$ws.Range("L132").Value = 8188.7145
$ws.Range("M132").Value = -2292.7502
$ws.Range("N132").Value = -13248.7145
$ws.Range("H133").Value = 48000
$ws.Range("J133").Value = 48000
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -53060
$ws.Range("H136").Value = 1886.875
$ws.Range("I136").Value = 1277.85
$ws.Range("K136").Value = 3833.55
$ws.Range("M136").Value = -1283.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2527.2354
$ws.Range("I134").Value = 2453.5
$ws.Range("J134").Value = 2871.3333
$ws.Range("K134").Value = 7360.5
$ws.Range("L134").Value = 8613.999899999999
$ws.Range("M134").Value = -4825.5
$ws.Range("N134").Value = -13683.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 14962.643
$ws.Range("I33").Value = 1206
$ws.Range("J33").Value = 28719.285
$ws.Range("K33").Value = 1206
$ws.Range("L33").Value = 28719.285
$ws.Range("M33").Value = -827
$ws.Range("N33").Value = -29477.285
$ws.Range("H58").Value = 1119.4878
$ws.Range("I58").Value = 946.40625
$ws.Range("J58").Value = 1734.8889
$ws.Range("K58").Value = 946.40625
$ws.Range("L58").Value = 1734.8889
$ws.Range("M58").Value = -743.40625
$ws.Range("N58").Value = -2140.8889
$ws.Range("H99").Value = 2533.3333
$ws.Range("I99").Value = 2114.2856
$ws.Range("K99").Value = 2114.2856
$ws.Range("M99").Value = -616.2856000000002
$ws.Range("H126").Value = 2533.3333
$ws.Range("I126").Value = 2114.2856
$ws.Range("K126").Value = 6342.8568
$ws.Range("M126").Value = -3872.8568
$ws.Range("H132").Value = 1792.6666
$ws.Range("I132").Value = 1491.0385
$ws.Range("J132").Value = 3753.25
$ws.Range("K132").Value = 4473.1155
$ws.Range("L132").Value = 11259.75
$ws.Range("M132").Value = -1943.1155
$ws.Range("N132").Value = -16319.75
$ws.Range("H134").Value = 1983.7742
$ws.Range("I134").Value = 1203.4584
$ws.Range("J134").Value = 4659.143
$ws.Range("K134").Value = 3610.3752
$ws.Range("L134").Value = 13977.429
$ws.Range("M134").Value = -1075.3752
$ws.Range("N134").Value = -19047.429
$ws.Range("H136").Value = 1119.4878
$ws.Range("I136").Value = 946.40625
$ws.Range("J136").Value = 1734.8889
$ws.Range("K136").Value = 2839.21875
$ws.Range("L136").Value = 5204.6667
$ws.Range("M136").Value = -289.21875
$ws.Range("N136").Value = -10304.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 611.25
$ws.Range("I92").Value = 598.3333
$ws.Range("J92").Value = 650
$ws.Range("K92").Value = 1794.9999
$ws.Range("L92").Value = 1950
$ws.Range("M92").Value = -546.9999
$ws.Range("N92").Value = -4446
$ws.Range("H94").Value = 3003.8333
$ws.Range("I94").Value = 1011.5
$ws.Range("K94").Value = 3034.5
$ws.Range("M94").Value = -2358.5
$ws.Range("H132").Value = 674845.9
$ws.Range("I132").Value = 1307.1818
$ws.Range("J132").Value = 2527077.2
$ws.Range("K132").Value = 11764.6362
$ws.Range("L132").Value = 22743694.8
$ws.Range("M132").Value = -9234.636200000001
$ws.Range("N132").Value = -22748754.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12029.857
$ws.Range("I70").Value = 13116.667
$ws.Range("J70").Value = 5509
$ws.Range("K70").Value = 13116.667
$ws.Range("L70").Value = 5509
$ws.Range("M70").Value = -12846.667
$ws.Range("N70").Value = -6049
$ws.Range("H73").Value = 12029.857
$ws.Range("I73").Value = 13116.667
$ws.Range("J73").Value = 5509
$ws.Range("K73").Value = 13116.667
$ws.Range("L73").Value = 5509
$ws.Range("M73").Value = -12180.667
$ws.Range("N73").Value = -7381
$ws.Range("H80").Value = 3015.1428
$ws.Range("I80").Value = 2820
$ws.Range("J80").Value = 3503
$ws.Range("K80").Value = 2820
$ws.Range("L80").Value = 3503
$ws.Range("M80").Value = -1822
$ws.Range("N80").Value = -5499
$ws.Range("H83").Value = 3015.1428
$ws.Range("I83").Value = 2820
$ws.Range("J83").Value = 3503
$ws.Range("K83").Value = 14100
$ws.Range("L83").Value = 17515
$ws.Range("M83").Value = -9108
$ws.Range("N83").Value = -27499
$ws.Range("H113").Value = 6751.15
$ws.Range("I113").Value = 1264.5454
$ws.Range("J113").Value = 13457
$ws.Range("K113").Value = 1264.5454
$ws.Range("L113").Value = 13457
$ws.Range("M113").Value = 905.4546
$ws.Range("N113").Value = -17797
$ws.Range("H132").Value = 2383.804
$ws.Range("I132").Value = 2015.3422
$ws.Range("J132").Value = 3460.8462
$ws.Range("K132").Value = 6046.0266
$ws.Range("L132").Value = 10382.5386
$ws.Range("M132").Value = -3516.0266
$ws.Range("N132").Value = -15442.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 55007.5
$ws.Range("J30").Value = 55007.5
$ws.Range("L30").Value = 55007.5
$ws.Range("N30").Value = -55221.5
$ws.Range("H81").Value = 1368.5714
$ws.Range("J81").Value = 1860
$ws.Range("L81").Value = 3720
$ws.Range("N81").Value = -5842
$ws.Range("H84").Value = 1368.5714
$ws.Range("J84").Value = 1860
$ws.Range("L84").Value = 18600
$ws.Range("N84").Value = -29208
